# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the regenerated data output.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 449
$ws1.Range("F4").Value = 1450
$ws1.Range("F5").Value = 141
$ws1.Range("F6").Value = 1731
$ws1.Range("F8").Value = 141
$ws1.Range("F9").Value = 652
$ws1.Range("F12").Value = 550
$ws1.Range("F15").Value = 139
$ws1.Range("F18").Value = 66
$ws1.Range("F19").Value = 101
$ws1.Range("F20").Value = 4530
$ws1.Range("F21").Value = 38
$ws1.Range("F22").Value = 809
$ws1.Range("F24").Value = 2158
$ws1.Range("F26").Value = 4
$ws1.Range("F27").Value = 2025

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 449
$ws4.Range("F4").Value = 1450
$ws4.Range("F5").Value = 141
$ws4.Range("F6").Value = 1731
$ws4.Range("F8").Value = 141
$ws4.Range("F9").Value = 652
$ws4.Range("F12").Value = 550
$ws4.Range("F15").Value = 139
$ws4.Range("F18").Value = 66
$ws4.Range("F19").Value = 101
$ws4.Range("F20").Value = 4530
$ws4.Range("F22").Value = 38
$ws4.Range("F24").Value = 809
$ws4.Range("F26").Value = 2158
$ws4.Range("F28").Value = 4
$ws4.Range("F29").Value = 2025
